# "Fruta / hortaliza, semanal"
#
# Weekly refresh of the "Vega Monumental Concepción - Papa" sheet: a new
# survey row (newest date) is inserted at the top of the data block (row 69,
# right after the first 67 data rows which stay untouched), which pushes
# every following record down by one row. The previously-last record simply
# ends up one row lower, at the new bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 69 - this shifts rows 69..152 down to 70..153,
# preserving all of their data/styles, and grows the used range to R153.
$ws.Rows.Item(69).Insert()

# Populate the newly inserted row with the new survey record.
$ws.Cells.Item(69, 1).Value  = 11
$ws.Cells.Item(69, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(69, 3).Value  = "Bíobío"
$ws.Cells.Item(69, 4).Value  = 44546
$ws.Cells.Item(69, 5).Value  = 8
$ws.Cells.Item(69, 6).Value  = 100114001
$ws.Cells.Item(69, 7).Value  = "Papa"
$ws.Cells.Item(69, 8).Value  = "Asterix"
$ws.Cells.Item(69, 9).Value  = "1a (nueva lavada)"
$ws.Cells.Item(69, 10).Value = 450
$ws.Cells.Item(69, 11).Value = 10000
$ws.Cells.Item(69, 12).Value = 11000
$ws.Cells.Item(69, 13).Value = 10444
$ws.Cells.Item(69, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(69, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(69, 16).Value = 418
$ws.Cells.Item(69, 17).Value = 25
$ws.Cells.Item(69, 18).Value = "Hortaliza"
